$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34
$prev = $row - 1

# Copy formatting from the previous data row (33) so the new row matches
# the existing style pattern (bold/bordered index column, date-formatted
# match-date column) without introducing new style entries.
$ws.Range("A$prev`:V$prev").Copy()
$ws.Range("A$row`:V$row").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 33
$ws.Cells.Item($row, 2).Value = "kuwait"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45233.70833333334
$ws.Cells.Item($row, 6).Value = "Al Shabab"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Al Jahra"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 2.33
$ws.Cells.Item($row, 11).Value = "02/11/2023 08:42"
$ws.Cells.Item($row, 12).Value = 2.16
$ws.Cells.Item($row, 13).Value = "03/11/2023 16:57"
$ws.Cells.Item($row, 14).Value = 3.18
$ws.Cells.Item($row, 15).Value = "02/11/2023 08:42"
$ws.Cells.Item($row, 16).Value = 3.34
$ws.Cells.Item($row, 17).Value = "03/11/2023 16:57"
$ws.Cells.Item($row, 18).Value = 2.6
$ws.Cells.Item($row, 19).Value = "02/11/2023 08:42"
$ws.Cells.Item($row, 20).Value = 3.11
$ws.Cells.Item($row, 21).Value = "03/11/2023 16:57"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-shabab-al-jahra/ddN5Ioe5/"
